# Applies the "456a3b4" data refresh:
#   - "展览"   sheet: F2 423->422, F3 2821->2834, F4 130->131, append row 5
#   - "全部类型" sheet: F2 423->422, F7 2821->2834, F8 130->131, append row 10
#
# Notes:
#  - New-row text cells that look like ISO dates (e.g. "2024-12-07") would be
#    auto-coerced to date serials by a normal typed-Value assignment, so that
#    cell is entered as a `="literal"` formula and then flattened back to a
#    plain value via Copy/PasteSpecial(xlPasteValuesAndNumberFormats) - this
#    avoids minting a stray number-format style too.
#  - The leading "index" column (A) uses a bold/centered/bordered style in
#    every existing data row; Range.Copy() onto the new cell (before writing
#    its value) reuses that existing style instead of minting a new one.
#  - Named parameters on custom functions aren't reliable here, so the
#    helper below uses positional parameters.

$wb = $excel.ActiveWorkbook

function Set-RowData($ws, $row, $index, $date, $name, $loc, $timeRange, $wantCount, $minPrice, $link, $cover) {
    # A: numeric index, styled like the other rows in column A.
    $aCell = $ws.Cells.Item($row, 1)
    $ws.Cells.Item($row - 1, 1).Copy($aCell)
    $aCell.Value = $index

    # B: date-like text - must stay plain text, not become a date serial.
    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Formula = '="' + $date + '"'
    $bCell.Copy()
    $bCell.PasteSpecial(-4163) | Out-Null

    $ws.Cells.Item($row, 3).Value = $name
    $ws.Cells.Item($row, 4).Value = $loc
    $ws.Cells.Item($row, 5).Value = $timeRange
    $ws.Cells.Item($row, 6).Value = $wantCount
    $ws.Cells.Item($row, 7).Value = $minPrice
    $ws.Cells.Item($row, 8).Value = $link
    $ws.Cells.Item($row, 9).Value = $cover
}

$newName = "合肥·心动恋章·冬日序国乙&代号鸢同人only"
$newLoc = "上海路与迎淮路交口向东200米(云峯中心一楼) 费加罗宴会艺术中心(省府店)"
$newTimeRange = "2024.12.07 12:00-12.07 21:00"
$newLink = "https://show.bilibili.com/platform/detail.html?id=93319"
$newCover = "//i0.hdslb.com/bfs/openplatform/202409/KtMLL8ZO1727684987784.jpeg"

# ---- "展览" sheet ----
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 422
$wsExpo.Range("F3").Value = 2834
$wsExpo.Range("F4").Value = 131

Set-RowData $wsExpo 5 4 "2024-12-07" $newName $newLoc $newTimeRange 8 50 $newLink $newCover

# ---- "全部类型" sheet ----
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 422
$wsAll.Range("F7").Value = 2834
$wsAll.Range("F8").Value = 131

Set-RowData $wsAll 10 9 "2024-12-07" $newName $newLoc $newTimeRange 8 50 $newLink $newCover

Write-Output "applied"
